$wb = $excel.ActiveWorkbook

# Sheet '展览' (index 1) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 250
$ws1.Cells.Item(3, 6).Value = 239
$ws1.Cells.Item(4, 6).Value = 260
$ws1.Cells.Item(5, 6).Value = 2861
$ws1.Cells.Item(7, 6).Value = 121
$ws1.Cells.Item(8, 6).Value = 2216
$ws1.Cells.Item(9, 6).Value = 346
$ws1.Cells.Item(11, 6).Value = 431
$ws1.Cells.Item(12, 6).Value = 81
$ws1.Cells.Item(13, 6).Value = 2544
$ws1.Cells.Item(15, 6).Value = 1333
$ws1.Cells.Item(16, 6).Value = 4681
$ws1.Cells.Item(17, 6).Value = 7
$ws1.Cells.Item(18, 6).Value = 5068
$ws1.Cells.Item(19, 6).Value = 1626
$ws1.Cells.Item(20, 6).Value = 2855
$ws1.Cells.Item(21, 6).Value = 3250
$ws1.Cells.Item(22, 6).Value = 162
$ws1.Cells.Item(23, 6).Value = 1547
$ws1.Cells.Item(24, 6).Value = 254
$ws1.Cells.Item(25, 6).Value = 834
$ws1.Cells.Item(26, 6).Value = 105
$ws1.Cells.Item(27, 6).Value = 288
$ws1.Cells.Item(28, 6).Value = 975
$ws1.Cells.Item(29, 6).Value = 1820
$ws1.Cells.Item(30, 6).Value = 116
$ws1.Cells.Item(31, 6).Value = 276
$ws1.Cells.Item(32, 6).Value = 693
$ws1.Cells.Item(34, 6).Value = 331
$ws1.Cells.Item(35, 6).Value = 406

# Sheet '全部类型' (index 4) - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(7, 6).Value = 250
$ws4.Cells.Item(8, 6).Value = 239
$ws4.Cells.Item(10, 6).Value = 260
$ws4.Cells.Item(11, 6).Value = 2861
$ws4.Cells.Item(13, 6).Value = 2216
$ws4.Cells.Item(14, 6).Value = 346
$ws4.Cells.Item(18, 6).Value = 431
$ws4.Cells.Item(19, 6).Value = 81
$ws4.Cells.Item(21, 6).Value = 2544
$ws4.Cells.Item(22, 6).Value = 1333
$ws4.Cells.Item(26, 6).Value = 4681
$ws4.Cells.Item(27, 6).Value = 7
$ws4.Cells.Item(28, 6).Value = 5068
$ws4.Cells.Item(29, 6).Value = 1626
$ws4.Cells.Item(30, 6).Value = 2855
$ws4.Cells.Item(31, 6).Value = 3250
$ws4.Cells.Item(32, 6).Value = 162
$ws4.Cells.Item(35, 6).Value = 1547
$ws4.Cells.Item(37, 6).Value = 254
$ws4.Cells.Item(38, 6).Value = 834
$ws4.Cells.Item(39, 6).Value = 105
$ws4.Cells.Item(40, 6).Value = 288
$ws4.Cells.Item(41, 6).Value = 975
$ws4.Cells.Item(43, 6).Value = 1820
$ws4.Cells.Item(44, 6).Value = 116
$ws4.Cells.Item(45, 6).Value = 276
$ws4.Cells.Item(46, 6).Value = 693
$ws4.Cells.Item(48, 6).Value = 331
$ws4.Cells.Item(49, 6).Value = 406
